# Pooh Points: normal 20260207
# Applies the score/status update for the ARK@MSST game (clock 13:06 -> 11:29, 1st Half)
# plus the associated stat-line corrections on the Players sheet, a 3-row player
# reassignment (Undrafted bench rows 13-15), the OwnerTotals starter totals, and the
# "min" column width widening on the Players sheet.

$wb = $excel.ActiveWorkbook
$wsPlayers = $wb.Worksheets.Item("Players")
$wsOwners = $wb.Worksheets.Item("OwnerTotals")

# Widen column P ("min") from raw width 5 to raw width 6.
$wsPlayers.Range("P1").EntireColumn.ColumnWidth = 5.14

# --- Players sheet cell updates ---
$wsPlayers.Range("G2").Value = '11:29 - 1st Half'
$wsPlayers.Range("O2").Value = 2
$wsPlayers.Range("P2").Value = 8
$wsPlayers.Range("G3").Value = '11:29 - 1st Half'
$wsPlayers.Range("H3").Value = 15
$wsPlayers.Range("I3").Value = 10
$wsPlayers.Range("J3").Value = 3
$wsPlayers.Range("L3").Value = 1
$wsPlayers.Range("M3").Value = 3
$wsPlayers.Range("P3").Value = 19
$wsPlayers.Range("Q3").Value = 4
$wsPlayers.Range("R3").Value = 6
$wsPlayers.Range("T3").Value = 3
$wsPlayers.Range("V3").Value = 1
$wsPlayers.Range("G4").Value = '11:29 - 1st Half'
$wsPlayers.Range("H4").Value = 3
$wsPlayers.Range("I4").Value = 7
$wsPlayers.Range("L4").Value = 1
$wsPlayers.Range("P4").Value = 19
$wsPlayers.Range("Q4").Value = 3
$wsPlayers.Range("R4").Value = 10
$wsPlayers.Range("S4").Value = 1
$wsPlayers.Range("T4").Value = 4
$wsPlayers.Range("G5").Value = '11:29 - 1st Half'
$wsPlayers.Range("H5").Value = -3
$wsPlayers.Range("K5").Value = 2
$wsPlayers.Range("N5").Value = 2
$wsPlayers.Range("P5").Value = 12
$wsPlayers.Range("R5").Value = 3
$wsPlayers.Range("G6").Value = '11:29 - 1st Half'
$wsPlayers.Range("O6").Value = 2
$wsPlayers.Range("P6").Value = 9
$wsPlayers.Range("G7").Value = '11:29 - 1st Half'
$wsPlayers.Range("H7").Value = 3
$wsPlayers.Range("I7").Value = 1
$wsPlayers.Range("J7").Value = 4
$wsPlayers.Range("N7").Value = 1
$wsPlayers.Range("P7").Value = 10
$wsPlayers.Range("U7").Value = 1
$wsPlayers.Range("V7").Value = 2
$wsPlayers.Range("G8").Value = '11:29 - 1st Half'
$wsPlayers.Range("H8").Value = 12
$wsPlayers.Range("I8").Value = 9
$wsPlayers.Range("J8").Value = 4
$wsPlayers.Range("K8").Value = 5
$wsPlayers.Range("N8").Value = 1
$wsPlayers.Range("P8").Value = 19
$wsPlayers.Range("Q8").Value = 4
$wsPlayers.Range("R8").Value = 9
$wsPlayers.Range("T8").Value = 1
$wsPlayers.Range("U8").Value = 1
$wsPlayers.Range("V8").Value = 1
$wsPlayers.Range("G9").Value = '11:29 - 1st Half'
$wsPlayers.Range("H9").Value = 7
$wsPlayers.Range("I9").Value = 10
$wsPlayers.Range("N9").Value = 2
$wsPlayers.Range("P9").Value = 15
$wsPlayers.Range("Q9").Value = 4
$wsPlayers.Range("R9").Value = 5
$wsPlayers.Range("U9").Value = 2
$wsPlayers.Range("V9").Value = 3
$wsPlayers.Range("G10").Value = '11:29 - 1st Half'
$wsPlayers.Range("H10").Value = 5
$wsPlayers.Range("I10").Value = 2
$wsPlayers.Range("J10").Value = 4
$wsPlayers.Range("P10").Value = 14
$wsPlayers.Range("Q10").Value = 1
$wsPlayers.Range("R10").Value = 3
$wsPlayers.Range("G11").Value = '11:29 - 1st Half'
$wsPlayers.Range("H11").Value = 11
$wsPlayers.Range("I11").Value = 9
$wsPlayers.Range("J11").Value = 5
$wsPlayers.Range("L11").Value = 1
$wsPlayers.Range("N11").Value = 1
$wsPlayers.Range("P11").Value = 15
$wsPlayers.Range("Q11").Value = 4
$wsPlayers.Range("R11").Value = 7
$wsPlayers.Range("S11").Value = 1
$wsPlayers.Range("T11").Value = 1
$wsPlayers.Range("G12").Value = '11:29 - 1st Half'
$wsPlayers.Range("H12").Value = 7
$wsPlayers.Range("J12").Value = 2
$wsPlayers.Range("L12").Value = 2
$wsPlayers.Range("M12").Value = 1
$wsPlayers.Range("N12").Value = 1
$wsPlayers.Range("O12").Value = 1
$wsPlayers.Range("P12").Value = 10
$wsPlayers.Range("D13").Value = 'Jamarion Davis-Fleming'
$wsPlayers.Range("E13").Value = 'MSST'
$wsPlayers.Range("G13").Value = '11:29 - 1st Half'
$wsPlayers.Range("H13").Value = 6
$wsPlayers.Range("I13").Value = 4
$wsPlayers.Range("J13").Value = 3
$wsPlayers.Range("O13").Value = 2
$wsPlayers.Range("P13").Value = 10
$wsPlayers.Range("Q13").Value = 2
$wsPlayers.Range("R13").Value = 3
$wsPlayers.Range("D14").Value = 'Isaiah Sealy'
$wsPlayers.Range("E14").Value = 'ARK'
$wsPlayers.Range("G14").Value = '11:29 - 1st Half'
$wsPlayers.Range("H14").Value = 4
$wsPlayers.Range("I14").Value = 2
$wsPlayers.Range("J14").Value = 1
$wsPlayers.Range("K14").Value = 0
$wsPlayers.Range("M14").Value = 3
$wsPlayers.Range("N14").Value = 1
$wsPlayers.Range("O14").Value = 3
$wsPlayers.Range("P14").Value = 10
$wsPlayers.Range("Q14").Value = 1
$wsPlayers.Range("R14").Value = 2
$wsPlayers.Range("D15").Value = 'Ja''Borri McGhee'
$wsPlayers.Range("G15").Value = '11:29 - 1st Half'
$wsPlayers.Range("H15").Value = 4
$wsPlayers.Range("I15").Value = 2
$wsPlayers.Range("J15").Value = 2
$wsPlayers.Range("K15").Value = 2
$wsPlayers.Range("O15").Value = 0
$wsPlayers.Range("P15").Value = 10
$wsPlayers.Range("R15").Value = 2
$wsPlayers.Range("T15").Value = 1
$wsPlayers.Range("U15").Value = 2
$wsPlayers.Range("V15").Value = 2
$wsPlayers.Range("G16").Value = '11:29 - 1st Half'
$wsPlayers.Range("H16").Value = -1
$wsPlayers.Range("J16").Value = 2
$wsPlayers.Range("N16").Value = 1
$wsPlayers.Range("P16").Value = 5
$wsPlayers.Range("R16").Value = 2
$wsPlayers.Range("T16").Value = 1
$wsPlayers.Range("G17").Value = '11:29 - 1st Half'

# --- OwnerTotals sheet cell updates ---
$wsOwners.Range("B2").Value = 15
$wsOwners.Range("B3").Value = 12
